$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 500
$ws.Range("I13").Value = 500
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 500
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -331
$ws.Range("N13").ClearContents()
$ws.Range("H19").Value = 3487.2
$ws.Range("I19").Value = 2998.3333
$ws.Range("K19").Value = 2998.3333
$ws.Range("M19").Value = -2823.3333
$ws.Range("H31").Value = 8000
$ws.Range("I31").Value = 8000
$ws.Range("K31").Value = 24000
$ws.Range("M31").Value = -23770
$ws.Range("H106").Value = 2151.25
$ws.Range("I106").Value = 2148
$ws.Range("K106").Value = 2148
$ws.Range("M106").Value = -1517
$ws.Range("H129").Value = 4167849.2
$ws.Range("I129").Value = 1314.2222
$ws.Range("K129").Value = 3942.6666
$ws.Range("M129").Value = 1057.3334
$ws.Range("H132").Value = 1472.138
$ws.Range("I132").Value = 1256.2273
$ws.Range("K132").Value = 3768.6819
$ws.Range("M132").Value = -1238.6819
$ws.Range("H139").Value = 150000
$ws.Range("J139").Value = 150000
$ws.Range("L139").Value = 150000
$ws.Range("N139").Value = -160280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2293.0212
$ws.Range("I32").Value = 2210.75
$ws.Range("J32").Value = 3499.6667
$ws.Range("K32").Value = 2210.75
$ws.Range("L32").Value = 3499.6667
$ws.Range("M32").Value = -1923.75
$ws.Range("N32").Value = -4073.6667
$ws.Range("H45").Value = 114496.78
$ws.Range("I45").Value = 155937.84
$ws.Range("K45").Value = 155937.84
$ws.Range("M45").Value = -155560.84
$ws.Range("H61").Value = 8739.291999999999
$ws.Range("I61").Value = 7844.905
$ws.Range("K61").Value = 7844.905
$ws.Range("M61").Value = -7632.905
$ws.Range("H97").Value = 826.6667
$ws.Range("I97").Value = 932
$ws.Range("J97").Value = 300
$ws.Range("K97").Value = 932
$ws.Range("L97").Value = 300
$ws.Range("M97").Value = -436
$ws.Range("N97").Value = -1292
$ws.Range("H122").Value = 2811.5
$ws.Range("I122").Value = 2525.9473
$ws.Range("J122").Value = 3304.7273
$ws.Range("K122").Value = 7577.841899999999
$ws.Range("L122").Value = 9914.1819
$ws.Range("M122").Value = -5127.841899999999
$ws.Range("N122").Value = -14814.1819
$ws.Range("H132").Value = 7007.091
$ws.Range("I132").Value = 5996.357
$ws.Range("K132").Value = 17989.071
$ws.Range("M132").Value = -15459.071
$ws.Range("H136").Value = 8739.291999999999
$ws.Range("I136").Value = 7844.905
$ws.Range("K136").Value = 23534.715
$ws.Range("M136").Value = -20984.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 273.33334
$ws.Range("I105").Value = 278
$ws.Range("J105").Value = 250
$ws.Range("K105").Value = 278
$ws.Range("L105").Value = 250
$ws.Range("M105").Value = 1469
$ws.Range("N105").Value = -3744
$ws.Range("H134").Value = 3408.2368
$ws.Range("I134").Value = 3488.1892
$ws.Range("J134").Value = 450
$ws.Range("K134").Value = 10464.5676
$ws.Range("L134").Value = 1350
$ws.Range("M134").Value = -7929.567599999998
$ws.Range("N134").Value = -6420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 29997.5
$ws.Range("J53").Value = 29997.5
$ws.Range("L53").Value = 29997.5
$ws.Range("N53").Value = -31211.5
$ws.Range("H99").Value = 4512.3335
$ws.Range("I99").Value = 4044.8
$ws.Range("J99").Value = 6850
$ws.Range("K99").Value = 4044.8
$ws.Range("L99").Value = 6850
$ws.Range("M99").Value = -2546.8
$ws.Range("N99").Value = -9846
$ws.Range("H126").Value = 4512.3335
$ws.Range("I126").Value = 4044.8
$ws.Range("J126").Value = 6850
$ws.Range("K126").Value = 12134.4
$ws.Range("L126").Value = 20550
$ws.Range("M126").Value = -9664.400000000001
$ws.Range("N126").Value = -25490
$ws.Range("H132").Value = 3822.5715
$ws.Range("I132").Value = 3550.4285
$ws.Range("K132").Value = 10651.2855
$ws.Range("M132").Value = -8121.2855
$ws.Range("H134").Value = 7238.6665
$ws.Range("I134").Value = 6281.125
$ws.Range("J134").Value = 10302.8
$ws.Range("K134").Value = 18843.375
$ws.Range("L134").Value = 30908.4
$ws.Range("M134").Value = -16308.375
$ws.Range("N134").Value = -35978.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 19997
$ws.Range("I125").Value = 19988
$ws.Range("K125").Value = 59964
$ws.Range("M125").Value = -55044
$ws.Range("H126").Value = 11499.833
$ws.Range("I126").Value = 7999.6665
$ws.Range("K126").Value = 23998.9995
$ws.Range("M126").Value = -19058.9995
$ws.Range("H128").Value = 1017071.3
$ws.Range("I128").Value = 1017071.3
$ws.Range("K128").Value = 3051213.9
$ws.Range("M128").Value = -3046233.9
$ws.Range("H129").Value = 20261776
$ws.Range("J129").Value = 1263112.4
$ws.Range("L129").Value = 3789337.2
$ws.Range("N129").Value = -3799337.2
$ws.Range("H130").Value = 7089.4
$ws.Range("I130").Value = 1815.6666
$ws.Range("K130").Value = 5446.9998
$ws.Range("M130").Value = -426.9997999999996
$ws.Range("H137").Value = 2271.6667
$ws.Range("J137").Value = 5000
$ws.Range("L137").Value = 15000
$ws.Range("N137").Value = -25200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4734.067
$ws.Range("I122").Value = 6040.875
$ws.Range("J122").Value = 3240.5715
$ws.Range("K122").Value = 18122.625
$ws.Range("L122").Value = 9721.7145
$ws.Range("M122").Value = -15672.625
$ws.Range("N122").Value = -14621.7145
$ws.Range("H126").Value = 3339398.8
$ws.Range("I126").Value = 6670299.5
$ws.Range("K126").Value = 20010898.5
$ws.Range("M126").Value = -20008428.5
$ws.Range("H132").Value = 2096.926
$ws.Range("I132").Value = 1912.72
$ws.Range("K132").Value = 5738.16
$ws.Range("M132").Value = -3208.16

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 13782.759
$ws.Range("J46").Value = 16751.715
$ws.Range("L46").Value = 16751.715
$ws.Range("N46").Value = -17127.715
$ws.Range("H132").Value = 8029.9546
$ws.Range("I132").Value = 7804.3945
$ws.Range("K132").Value = 23413.1835
$ws.Range("M132").Value = -20883.1835

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3693.2285
$ws.Range("I122").Value = 3465.2222
$ws.Range("J122").Value = 4462.75
$ws.Range("K122").Value = 10395.6666
$ws.Range("L122").Value = 13388.25
$ws.Range("M122").Value = -7945.6666
$ws.Range("N122").Value = -18288.25
$ws.Range("H132").Value = 7549.6
$ws.Range("I132").Value = 7070.857
$ws.Range("K132").Value = 21212.571
$ws.Range("M132").Value = -18682.571
$ws.Range("H135").Value = 135959.17
$ws.Range("J135").Value = 147021
$ws.Range("L135").Value = 147021
$ws.Range("N135").Value = -157161
$ws.Range("H136").Value = 3282.4146
$ws.Range("I136").Value = 2073.9143
$ws.Range("J136").Value = 10332
$ws.Range("K136").Value = 6221.742899999999
$ws.Range("L136").Value = 30996
$ws.Range("M136").Value = -3671.742899999999
$ws.Range("N136").Value = -36096
